$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.990137338638306
$ws.Range("B1").Value = 5.348728179931641
$ws.Range("C1").Value = 2.293470859527588
$ws.Range("D1").Value = 1.367118835449219
$ws.Range("E1").Value = 1.372858762741089
